# Generate Report for Handoff
# Adds two new tracked files (4bd41fab-... and 5aeaf081-...) to the
# localization-status workbook: one new row each on the "Overview",
# "zh-cn" and "de-de" sheets, with matching hyperlinks.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet - rows 4 and 5
# ---------------------------------------------------------------------

$overview.Range("A4").Value = "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md"
$overview.Hyperlinks.Add($overview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/e2e/4bd41fab-3676-40ef-a41b-179d0ebe3daa.md", "", "", "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md") | Out-Null
$overview.Range("A4").Style = "HyperLink"
$overview.Range("B4").Value = "Ready for handoff"
$overview.Range("C4").Value = "Ready for handoff"
$overview.Range("D4").Value = "2016-40-18 05:40:40"

$overview.Range("A5").Value = "5aeaf081-63b8-434f-8808-e8aabd8a537e.md"
$overview.Hyperlinks.Add($overview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/e2e/5aeaf081-63b8-434f-8808-e8aabd8a537e.md", "", "", "5aeaf081-63b8-434f-8808-e8aabd8a537e.md") | Out-Null
$overview.Range("A5").Style = "HyperLink"
$overview.Range("B5").Value = "Ready for handoff"
$overview.Range("C5").Value = "Ready for handoff"
$overview.Range("D5").Value = "2016-40-18 05:40:40"

# ---------------------------------------------------------------------
# zh-cn sheet - rows 4 and 5
# ---------------------------------------------------------------------

$zhcn.Range("A4").Value = "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md"
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/e2e/4bd41fab-3676-40ef-a41b-179d0ebe3daa.md", "", "", "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md") | Out-Null
$zhcn.Range("A4").Style = "HyperLink"

$zhcn.Range("B4").Value = ".md"
$zhcn.Hyperlinks.Add($zhcn.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/e2e/4bd41fab-3676-40ef-a41b-179d0ebe3daa.md", "", "", ".md") | Out-Null
$zhcn.Range("B4").Style = "HyperLink"

$zhcn.Range("C4").Value = "Ready for handoff"

$zhcn.Range("D4").Value = "4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.zh-cn.xlf", "", "", "4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.zh-cn.xlf") | Out-Null
$zhcn.Range("D4").Style = "HyperLink"

$zhcn.Range("E4").Value = "2016-03-18 05:40:37"
$zhcn.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zhcn.Range("H4").Value = "0001-01-01 00:00:00"
$zhcn.Range("I4").Value = "Include"

$zhcn.Range("A5").Value = "5aeaf081-63b8-434f-8808-e8aabd8a537e.md"
$zhcn.Hyperlinks.Add($zhcn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/e2e/5aeaf081-63b8-434f-8808-e8aabd8a537e.md", "", "", "5aeaf081-63b8-434f-8808-e8aabd8a537e.md") | Out-Null
$zhcn.Range("A5").Style = "HyperLink"

$zhcn.Range("B5").Value = ".md"
$zhcn.Hyperlinks.Add($zhcn.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/e2e/5aeaf081-63b8-434f-8808-e8aabd8a537e.md", "", "", ".md") | Out-Null
$zhcn.Range("B5").Style = "HyperLink"

$zhcn.Range("C5").Value = "Ready for handoff"

$zhcn.Range("D5").Value = "5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.zh-cn.xlf", "", "", "5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.zh-cn.xlf") | Out-Null
$zhcn.Range("D5").Style = "HyperLink"

$zhcn.Range("E5").Value = "2016-03-18 05:40:37"
$zhcn.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zhcn.Range("H5").Value = "0001-01-01 00:00:00"
$zhcn.Range("I5").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet - rows 4 and 5
# ---------------------------------------------------------------------

$dede.Range("A4").Value = "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md"
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/e2e/4bd41fab-3676-40ef-a41b-179d0ebe3daa.md", "", "", "4bd41fab-3676-40ef-a41b-179d0ebe3daa.md") | Out-Null
$dede.Range("A4").Style = "HyperLink"

$dede.Range("B4").Value = ".md"
$dede.Hyperlinks.Add($dede.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/e2e/4bd41fab-3676-40ef-a41b-179d0ebe3daa.md", "", "", ".md") | Out-Null
$dede.Range("B4").Style = "HyperLink"

$dede.Range("C4").Value = "Ready for handoff"

$dede.Range("D4").Value = "4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bd41fab3676240ef240a41b240179d0ebe3daa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.de-de.xlf", "", "", "4bd41fab-3676-40ef-a41b-179d0ebe3daa.251d509209f64a784b423e814e2e19c9719dd05b.de-de.xlf") | Out-Null
$dede.Range("D4").Style = "HyperLink"

$dede.Range("E4").Value = "2016-03-18 05:40:40"
$dede.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$dede.Range("H4").Value = "0001-01-01 00:00:00"
$dede.Range("I4").Value = "Include"

$dede.Range("A5").Value = "5aeaf081-63b8-434f-8808-e8aabd8a537e.md"
$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/e2e/5aeaf081-63b8-434f-8808-e8aabd8a537e.md", "", "", "5aeaf081-63b8-434f-8808-e8aabd8a537e.md") | Out-Null
$dede.Range("A5").Style = "HyperLink"

$dede.Range("B5").Value = ".md"
$dede.Hyperlinks.Add($dede.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/e2e/5aeaf081-63b8-434f-8808-e8aabd8a537e.md", "", "", ".md") | Out-Null
$dede.Range("B5").Style = "HyperLink"

$dede.Range("C5").Value = "Ready for handoff"

$dede.Range("D5").Value = "5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5aeaf0815bb8834f238808e8aabd8a537e240000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.de-de.xlf", "", "", "5aeaf081-63b8-434f-8808-e8aabd8a537e.f09b48f1c0acb1e4af994c3cfb41d43fe2ebebcc.de-de.xlf") | Out-Null
$dede.Range("D5").Style = "HyperLink"

$dede.Range("E5").Value = "2016-03-18 05:40:40"
$dede.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$dede.Range("H5").Value = "0001-01-01 00:00:00"
$dede.Range("I5").Value = "Include"
